$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("!dData")

# The "(KCC0) Oficina" category no longer appears in the data; the row that used
# to reference it now references the existing "(KCC ) Oficina" label instead.
$ws.Range("E153").Value = "(KCC ) Oficina"

# Fill in the previously-zeroed "Pedidos" (orders) counts for rows 144-158.
$pedidos = @{
    144 = 3
    145 = 153
    146 = 57
    147 = 17
    148 = 1734
    149 = 1194
    150 = 8
    151 = 643
    152 = 294
    153 = 21
    154 = 2335
    155 = 1485
    156 = 9
    157 = 665
    158 = 347
}

foreach ($row in $pedidos.Keys) {
    $ws.Range("G$row").Value = $pedidos[$row]
}

# Reset the active selection back to the top-left of the table.
$ws.Range("B3").Select()
